$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 15 ---
$ws.Range("A15").Value = 44685
$ws.Range("B15").Value = 0.36805555555555558
$ws.Range("C15").Value = 0.39583333333333331
$ws.Range("E15").Value = "Migration sur Icescrum"
$ws.Range("F15").Value = "Sprint semain1"

# --- Row 16 ---
$ws.Range("A16").Value = 44685
$ws.Range("B16").Value = 0.40972222222222227
$ws.Range("C16").Value = 0.4381944444444445
$ws.Range("E16").Value = "Sprint Semaine 2"

# --- Row 17 ---
$ws.Range("A17").Value = 44685
$ws.Range("B17").Value = 0.4381944444444445
$ws.Range("C17").Value = 0.46388888888888885
$ws.Range("E17").Value = "Use Case"
$ws.Range("F17").Value = "Regiment et Unité: définition des concernes"

# --- Row 18 ---
$ws.Range("A18").Value = 44685
$ws.Range("B18").Value = 0.46388888888888885
$ws.Range("C18").Value = 0.49027777777777781

# --- Row 19 ---
$ws.Range("A19").Value = 44685
$ws.Range("B19").Value = 0.49027777777777781
$ws.Range("C19").Value = 0.50347222222222221
$ws.Range("E19").Value = "Implémentation du Menu"
$ws.Range("F19").Value = "Implémentattion des la partie graphique selon maquette"

# --- Row 20 ---
$ws.Range("A20").Value = 44685
$ws.Range("B20").Value = 0.56180555555555556
$ws.Range("C20").Value = 0.57777777777777783
$ws.Range("E20").Value = "Menu : Fonction"
$ws.Range("F20").Value = "Implémentation des fonctions liées au menu "

# --- Row 21 ---
$ws.Range("A21").Value = 44685
$ws.Range("B21").Value = 0.57777777777777783
$ws.Range("C21").Value = 0.61736111111111114
$ws.Range("E21").Value = "Documentation de`nla partie Tests de validation"
$ws.Range("E21").WrapText = $true
$ws.Range("F21").Value = "Menu principale`nIntéraction régiment"

# --- Row 22 ---
$ws.Range("A22").Value = 44685
$ws.Range("B22").Value = 0.61736111111111114
$ws.Range("C22").Value = 0.65763888888888888
$ws.Range("E22").Value = "Documentation"

# --- View state: scroll + selection ---
$ws.Activate()
$ws.Range("F22").Select()
